$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 44: new task entry under "Błędy magazynu" category
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "Błędy magazynu"
$ws.Cells.Item(44, 3).Value = "Kasowanie wpisów podczas powtórnego importu"
$ws.Cells.Item(44, 4).Value = 0

# Row 45: new task entry under "KPI magazynu " category
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "KPI magazynu "
$ws.Cells.Item(45, 3).Value = "Dodać możliwość wprowadzania przecinka i zamiany go na kropkę w czasie zapisu"
$ws.Cells.Item(45, 4).Value = 0

# Apply the same cell formatting used by the preceding rows
$ws.Range("A44:A45").VerticalAlignment = -4160
$ws.Range("B44:B45").VerticalAlignment = -4160
$ws.Range("C44:C45").WrapText = $true
$ws.Range("D44:D45").VerticalAlignment = -4108

# Update the selection to reflect where the user left off editing
$ws.Range("D46").Select()
